# Update "想去人数" (F column) values for matching rows in both the
# "展览" sheet and the "全部类型" sheet, per the commit's regenerated data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId=1) - row => new value for column F
$sheet1Updates = @{
    4  = 413
    5  = 1409
    7  = 2144
    11 = 4781
    15 = 220
    17 = 163
    21 = 3703
    22 = 575
    23 = 608
    24 = 24
    27 = 112
    29 = 11
    34 = 839
    35 = 2312
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (sheetId=4) - row => new value for column F
$sheet4Updates = @{
    4  = 413
    5  = 1409
    7  = 2144
    11 = 4781
    15 = 220
    17 = 163
    21 = 3703
    22 = 575
    23 = 608
    24 = 24
    27 = 112
    29 = 11
    35 = 839
    36 = 2312
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
